# Apply cryptos list update (prices / % changes / one row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.325.22"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.32"
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.50"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5043"
$ws.Range("E7").Value = "  -1.94%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09865"
$ws.Range("E9").Value = "  +26.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.129"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.45"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.479"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.02"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.870.13"
$ws.Range("E14").Value = "  +3.78%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.398"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("E17").Value = "  +5.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.68"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06673"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.115"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.381.34"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.262"
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.548"
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.45"
$ws.Range("E27").Value = "  +4.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.088.15"
$ws.Range("E28").Value = "  +3.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.76"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.03"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1065"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.654"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.606"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06836"
$ws.Range("E35").Value = "  -4.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.460"
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02394"
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.035"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.52"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6306"
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6026"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.281"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.683"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.55"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  +4.08%  "
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.127"
$ws.Range("E51").Value = "  +5.31%  "
